$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 321/322, shifting the former rows 321-397 down to 323-399.
$ws.Range("A321:A322").EntireRow.Insert()

# New row 321
$ws.Cells.Item(321, 1).Value = 10
$ws.Cells.Item(321, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(321, 3).Value = "La Araucanía"
$ws.Cells.Item(321, 4).Value = 44722
$ws.Cells.Item(321, 5).Value = 9
$ws.Cells.Item(321, 6).Value = 100114014
$ws.Cells.Item(321, 7).Value = "Betarraga"
$ws.Cells.Item(321, 8).Value = "Sin especificar"
$ws.Cells.Item(321, 9).Value = "Primera"
$ws.Cells.Item(321, 10).Value = 40
$ws.Cells.Item(321, 11).Value = 9500
$ws.Cells.Item(321, 12).Value = 9500
$ws.Cells.Item(321, 13).Value = 9500
$ws.Cells.Item(321, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(321, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(321, 16).Value = 792
$ws.Cells.Item(321, 17).Value = 12
$ws.Cells.Item(321, 18).Value = "Hortaliza"

# New row 322
$ws.Cells.Item(322, 1).Value = 10
$ws.Cells.Item(322, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(322, 3).Value = "La Araucanía"
$ws.Cells.Item(322, 4).Value = 44722
$ws.Cells.Item(322, 5).Value = 9
$ws.Cells.Item(322, 6).Value = 100114014
$ws.Cells.Item(322, 7).Value = "Betarraga"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 50
$ws.Cells.Item(322, 11).Value = 9500
$ws.Cells.Item(322, 12).Value = 9500
$ws.Cells.Item(322, 13).Value = 9500
$ws.Cells.Item(322, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(322, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(322, 16).Value = 380
$ws.Cells.Item(322, 17).Value = 25
$ws.Cells.Item(322, 18).Value = "Hortaliza"

# Ensure the date cells keep the same date-time number format used by the rest of column D.
$ws.Range("D321:D322").NumberFormat = $ws.Range("D323").NumberFormat
